$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: Aris, 2024-02-14, 8:07 - 9:35
$ws.Range("A24").Value = "Aris"
$ws.Range("B24").Value = 45336
$ws.Range("B23").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("C24").Formula = "=8+7/60"
$ws.Range("D24").Formula = "=9+35/60"

# Row 25: Viki, 2024-02-14, 8:07 - 9:35
$ws.Range("A25").Value = "Viki"
$ws.Range("B25").Value = 45336
$ws.Range("B23").Copy()
$ws.Range("B25").PasteSpecial(-4122)
$ws.Range("C25").Formula = "=8+7/60"
$ws.Range("D25").Formula = "=9+35/60"

# Update the selection to F25 to match the saved workbook state
$ws.Range("F25").Select()
